# Add PF/1.0.2 to meta-sheet
# Appends a new row (row 3) to Sheet1 with the new release tag "PF/1.0.2"
# in column A and "X" markers in columns B-D (sit2/uat2/prod not yet
# promoted for this release).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "PF/1.0.2"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
